# Add a "Save" column (H) to the s_vals sheet, matching the existing
# header styling used by the other header cells (copy style from G1),
# and set the two data rows' Save values to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold font, border, centered alignment) from the
# existing "sum" header (G1) onto the new "Save" header (H1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
